# "added further wage analysis"
# Rename the single worksheet from "JobTitleAndWage" to "County" and let the
# sheet view settle back on its default scroll position (drop the stale
# topLeftCell="A23" the sheet had been scrolled to) as the workbook is
# revisited for further analysis on a per-county basis.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "County"

# Reset the view to the top-left of the sheet (clears the old scrolled
# position) while keeping the existing selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G58").Select()
